$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Cells.Item(1, 1).Value = "Row"
$ws.Cells.Item(1, 2).Value = "Prognose"
$ws.Cells.Item(1, 3).Value = "surveys"
$ws.Cells.Item(1, 4).Value = "production"
$ws.Cells.Item(1, 5).Value = "orders"
$ws.Cells.Item(1, 6).Value = "turnover"
$ws.Cells.Item(1, 7).Value = "financial"
$ws.Cells.Item(1, 8).Value = "labor market"
$ws.Cells.Item(1, 9).Value = "prices"
$ws.Cells.Item(1, 10).Value = "national accounts"
$ws.Cells.Item(1, 11).Value = "Revision"

# --- Column A dates (rows 2-12), forced to Text so they are not
#     auto-recognized as date serials (matches the source data, which
#     stores these as plain strings) ---
$ws.Range("A2:A12").NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-03-30"
$ws.Cells.Item(3, 1).Value = "2025-04-15"
$ws.Cells.Item(4, 1).Value = "2025-04-30"
$ws.Cells.Item(5, 1).Value = "2025-05-15"
$ws.Cells.Item(6, 1).Value = "2025-05-30"
$ws.Cells.Item(7, 1).Value = "2025-06-15"
$ws.Cells.Item(8, 1).Value = "2025-06-30"
$ws.Cells.Item(9, 1).Value = "2025-07-15"
$ws.Cells.Item(10, 1).Value = "2025-07-30"
$ws.Cells.Item(11, 1).Value = "2025-08-15"
$ws.Cells.Item(12, 1).Value = "2025-08-30"
$ws.Range("A2:A12").ClearFormats()

# --- Numeric data B2:K12 ---
$ws.Cells.Item(2, 2).Value = 0.35803172283809664
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(3, 2).Value = 0.32441921910040877
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = -0.020850342990136762
$ws.Cells.Item(3, 5).Value = -0.0050559081910599929
$ws.Cells.Item(3, 6).Value = -0.0007180090807009491
$ws.Cells.Item(3, 7).Value = 0.0058852254117251384
$ws.Cells.Item(3, 8).Value = 0.0014779881212820701
$ws.Cells.Item(3, 9).Value = -0.017210167516011763
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0.0028587105072143681
$ws.Cells.Item(4, 2).Value = 0.33871481319809849
$ws.Cells.Item(4, 3).Value = 0.037584943393469128
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0.00095419891471768535
$ws.Cells.Item(4, 6).Value = -0.00043385831710949076
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0.00071811316639798224
$ws.Cells.Item(4, 9).Value = -0.025238268683045323
$ws.Cells.Item(4, 10).Value = 0.00056520748527336855
$ws.Cells.Item(4, 11).Value = 0.00014525813798638154
$ws.Cells.Item(5, 2).Value = 0.33303837312012863
$ws.Cells.Item(5, 3).Value = 0.10837641500756059
$ws.Cells.Item(5, 4).Value = -0.065920123911742387
$ws.Cells.Item(5, 5).Value = -0.00044721766904077436
$ws.Cells.Item(5, 6).Value = -0.020528872837830722
$ws.Cells.Item(5, 7).Value = -0.015151931277146908
$ws.Cells.Item(5, 8).Value = -0.0018307581827844018
$ws.Cells.Item(5, 9).Value = -0.010051249233989583
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = -0.00012270197299568686
$ws.Cells.Item(6, 2).Value = 0.41859258792155474
$ws.Cells.Item(6, 3).Value = 0.082155157599970066
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = -0.0038495691470361486
$ws.Cells.Item(6, 6).Value = 0.0030919005347437244
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0.0016523425949918399
$ws.Cells.Item(6, 9).Value = 0.000032669063740347996
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 0.0024717141550162913
$ws.Cells.Item(7, 2).Value = 0.46291474502950281
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0.018725842152795668
$ws.Cells.Item(7, 5).Value = -0.0089477585912688239
$ws.Cells.Item(7, 6).Value = -0.0018569273602913811
$ws.Cells.Item(7, 7).Value = 0.018379647315619144
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0.0017570818416337217
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0.016264271749459769
$ws.Cells.Item(8, 2).Value = 0.12914450367230632
$ws.Cells.Item(8, 3).Value = -0.28759318209023621
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = -0.00036538784589185065
$ws.Cells.Item(8, 6).Value = 0.0029881072297135362
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0.0009024479510256157
$ws.Cells.Item(8, 9).Value = -0.04918421646556044
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = -0.00051801013624719339
$ws.Cells.Item(9, 2).Value = 0.072640882065478835
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = -0.088581791935331777
$ws.Cells.Item(9, 5).Value = 0.0037466307495547945
$ws.Cells.Item(9, 6).Value = 0.015889177425760421
$ws.Cells.Item(9, 7).Value = 0.0044318096337087678
$ws.Cells.Item(9, 8).Value = 0.00036593117511360591
$ws.Cells.Item(9, 9).Value = 0.006247936228029086
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0.0013966851163376426
$ws.Cells.Item(10, 2).Value = 0.225721808646204
$ws.Cells.Item(10, 3).Value = 0.18450539919652911
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = -0.00038670275959956787
$ws.Cells.Item(10, 6).Value = 0.0014007031260898748
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = -0.00036882825622278173
$ws.Cells.Item(10, 9).Value = 0.015606097729508142
$ws.Cells.Item(10, 10).Value = -0.045292896107787937
$ws.Cells.Item(10, 11).Value = -0.0023828463477916506
$ws.Cells.Item(11, 2).Value = 0.4426047535069747
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0.13526485866657129
$ws.Cells.Item(11, 5).Value = -0.0024106126553314292
$ws.Cells.Item(11, 6).Value = 0.0033732158567832213
$ws.Cells.Item(11, 7).Value = -0.0022834126738595154
$ws.Cells.Item(11, 8).Value = -0.0030294172396428841
$ws.Cells.Item(11, 9).Value = 0.11059977466786804
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = -0.024631461761618012
$ws.Cells.Item(12, 2).Value = 0.18632258454614162
$ws.Cells.Item(12, 3).Value = -0.19863544576429448
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = -0.00036603963362944538
$ws.Cells.Item(12, 6).Value = 0.00011085146520785409
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = -0.000090154193706881404
$ws.Cells.Item(12, 9).Value = -0.036034833141970035
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = -0.021266547692440074

# --- Column width tweaks for G (7) and H (8) to reflect new content widths ---
$ws.Columns.Item(7).ColumnWidth = 14.333333333333334
$ws.Columns.Item(8).ColumnWidth = 15.333333333333334
